# Auto-generated update of Leve profit/price columns (H:N) across all sheets
# per scheduled runner refresh of market data.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 251.14815  # H33 was 209.34146
$ws.Cells.Item(33, 9).Value = 208.55  # I33 was 178.27272
$ws.Cells.Item(33, 10).Value = 372.85715  # J33 was 337.5
$ws.Cells.Item(33, 11).Value = 208.55  # K33 was 178.27272
$ws.Cells.Item(33, 12).Value = 372.85715  # L33 was 337.5
$ws.Cells.Item(33, 13).Value = 20.44999999999999  # M33 was 50.72728000000001
$ws.Cells.Item(33, 14).Value = -830.85715  # N33 was -795.5
$ws.Cells.Item(98, 8).Value = 986.46155  # H98 was 0
$ws.Cells.Item(98, 9).Value = 993.6667  # I98 was 0
$ws.Cells.Item(98, 10).Value = 900  # J98 was 0
$ws.Cells.Item(98, 11).Value = 993.6667  # K98 was 0
$ws.Cells.Item(98, 12).Value = 900  # L98 was 0
$ws.Cells.Item(98, 13).Value = 504.3333  # M98 was None
$ws.Cells.Item(98, 14).Value = -3896  # N98 was None
$ws.Cells.Item(111, 8).Value = 2056.2727  # H111 was 2191.8096
$ws.Cells.Item(111, 9).Value = 1569.6  # I111 was 1802
$ws.Cells.Item(111, 10).Value = 2461.8333  # J111 was 2431.6924
$ws.Cells.Item(111, 11).Value = 4708.799999999999  # K111 was 5406
$ws.Cells.Item(111, 12).Value = 7385.499899999999  # L111 was 7295.0772
$ws.Cells.Item(111, 13).Value = -1641.799999999999  # M111 was -2339
$ws.Cells.Item(111, 14).Value = -13519.4999  # N111 was -13429.0772
$ws.Cells.Item(122, 8).Value = 986.46155  # H122 was 0
$ws.Cells.Item(122, 9).Value = 993.6667  # I122 was 0
$ws.Cells.Item(122, 10).Value = 900  # J122 was 0
$ws.Cells.Item(122, 11).Value = 2981.0001  # K122 was 0
$ws.Cells.Item(122, 12).Value = 2700  # L122 was 0
$ws.Cells.Item(122, 13).Value = -531.0001000000002  # M122 was None
$ws.Cells.Item(122, 14).Value = -7600  # N122 was None
$ws.Cells.Item(125, 8).Value = 2106.6667  # H125 was 3000
$ws.Cells.Item(125, 9).Value = 320  # I125 was 0
$ws.Cells.Item(125, 11).Value = 2880  # K125 was 0
$ws.Cells.Item(125, 13).Value = -420  # M125 was None
$ws.Cells.Item(132, 8).Value = 792408.0600000001  # H132 was 1198713.5
$ws.Cells.Item(132, 9).Value = 2012.5454  # I132 was 3578.7144
$ws.Cells.Item(132, 10).Value = 7002658.5  # J132 was 8170333.5
$ws.Cells.Item(132, 11).Value = 6037.6362  # K132 was 10736.1432
$ws.Cells.Item(132, 12).Value = 21007975.5  # L132 was 24511000.5
$ws.Cells.Item(132, 13).Value = -3507.6362  # M132 was -8206.143199999999
$ws.Cells.Item(132, 14).Value = -21013035.5  # N132 was -24516060.5
$ws.Cells.Item(137, 8).Value = 1317577.5  # H137 was 1564544.8
$ws.Cells.Item(137, 9).Value = 1887856.8  # I137 was 2326805
$ws.Cells.Item(137, 10).Value = 3455.652  # J137 was 3726.3333
$ws.Cells.Item(137, 11).Value = 5663570.4  # K137 was 6980415
$ws.Cells.Item(137, 12).Value = 10366.956  # L137 was 11178.9999
$ws.Cells.Item(137, 13).Value = -5661020.4  # M137 was -6977865
$ws.Cells.Item(137, 14).Value = -15466.956  # N137 was -16278.9999
$ws.Cells.Item(138, 8).Value = 4881439  # H138 was 3776747.2
$ws.Cells.Item(138, 9).Value = 5297  # I138 was 8598.5
$ws.Cells.Item(138, 10).Value = 5885351  # J138 was 4084351
$ws.Cells.Item(138, 11).Value = 15891  # K138 was 25795.5
$ws.Cells.Item(138, 12).Value = 17656053  # L138 was 12253053
$ws.Cells.Item(138, 13).Value = -10751  # M138 was -20655.5
$ws.Cells.Item(138, 14).Value = -17666333  # N138 was -12263333

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2151.9644  # H2 was 2072.4827
$ws.Cells.Item(2, 9).Value = 1890.7059  # I2 was 1779.579
$ws.Cells.Item(2, 10).Value = 2555.7273  # J2 was 2629
$ws.Cells.Item(2, 11).Value = 1890.7059  # K2 was 1779.579
$ws.Cells.Item(2, 12).Value = 2555.7273  # L2 was 2629
$ws.Cells.Item(2, 13).Value = -1777.7059  # M2 was -1666.579
$ws.Cells.Item(2, 14).Value = -2781.7273  # N2 was -2855
$ws.Cells.Item(74, 8).Value = 5410332.5  # H74 was 5410284.5
$ws.Cells.Item(74, 9).Value = 7269274  # I74 was 7599628
$ws.Cells.Item(74, 10).Value = 65875  # J74 was 58555.555
$ws.Cells.Item(74, 11).Value = 7269274  # K74 was 7599628
$ws.Cells.Item(74, 12).Value = 65875  # L74 was 58555.555
$ws.Cells.Item(74, 13).Value = -7268400  # M74 was -7598754
$ws.Cells.Item(74, 14).Value = -67623  # N74 was -60303.555
$ws.Cells.Item(77, 8).Value = 5410332.5  # H77 was 5410284.5
$ws.Cells.Item(77, 9).Value = 7269274  # I77 was 7599628
$ws.Cells.Item(77, 10).Value = 65875  # J77 was 58555.555
$ws.Cells.Item(77, 11).Value = 36346370  # K77 was 37998140
$ws.Cells.Item(77, 12).Value = 329375  # L77 was 292777.775
$ws.Cells.Item(77, 13).Value = -36342002  # M77 was -37993772
$ws.Cells.Item(77, 14).Value = -338111  # N77 was -301513.775
$ws.Cells.Item(102, 8).Value = 3866.6667  # H102 was 7938712
$ws.Cells.Item(102, 9).Value = 3685.7144  # I102 was 10206233
$ws.Cells.Item(102, 10).Value = 4500  # J102 was 2388
$ws.Cells.Item(102, 11).Value = 3685.7144  # K102 was 10206233
$ws.Cells.Item(102, 12).Value = 4500  # L102 was 2388
$ws.Cells.Item(102, 13).Value = -2063.7144  # M102 was -10204611
$ws.Cells.Item(102, 14).Value = -7744  # N102 was -5632
$ws.Cells.Item(116, 8).Value = 2151.9644  # H116 was 2072.4827
$ws.Cells.Item(116, 9).Value = 1890.7059  # I116 was 1779.579
$ws.Cells.Item(116, 10).Value = 2555.7273  # J116 was 2629
$ws.Cells.Item(116, 11).Value = 1890.7059  # K116 was 1779.579
$ws.Cells.Item(116, 12).Value = 2555.7273  # L116 was 2629
$ws.Cells.Item(116, 13).Value = 403.2941000000001  # M116 was 514.421
$ws.Cells.Item(116, 14).Value = -7143.7273  # N116 was -7217
$ws.Cells.Item(132, 8).Value = 32217.508  # H132 was 35839.418
$ws.Cells.Item(132, 9).Value = 20646.604  # I132 was 22803.25
$ws.Cells.Item(132, 10).Value = 70546.125  # J132 was 80534.86
$ws.Cells.Item(132, 11).Value = 61939.812  # K132 was 68409.75
$ws.Cells.Item(132, 12).Value = 211638.375  # L132 was 241604.58
$ws.Cells.Item(132, 13).Value = -59409.812  # M132 was -65879.75
$ws.Cells.Item(132, 14).Value = -216698.375  # N132 was -246664.58

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2151.9644  # H3 was 2072.4827
$ws.Cells.Item(3, 9).Value = 1890.7059  # I3 was 1779.579
$ws.Cells.Item(3, 10).Value = 2555.7273  # J3 was 2629
$ws.Cells.Item(3, 11).Value = 1890.7059  # K3 was 1779.579
$ws.Cells.Item(3, 12).Value = 2555.7273  # L3 was 2629
$ws.Cells.Item(3, 13).Value = -1776.7059  # M3 was -1665.579
$ws.Cells.Item(3, 14).Value = -2783.7273  # N3 was -2857
$ws.Cells.Item(7, 8).Value = 2900  # H7 was 14930
$ws.Cells.Item(7, 10).Value = 0  # J7 was 38990
$ws.Cells.Item(7, 12).Value = 0  # L7 was 38990
$ws.Cells.Item(7, 14).ClearContents()  # N7 was -39216

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 234371.34  # H31 was 243730.3
$ws.Cells.Item(31, 9).Value = 32487.812  # I31 was 34627.168
$ws.Cells.Item(31, 11).Value = 32487.812  # K31 was 34627.168
$ws.Cells.Item(31, 13).Value = -32192.812  # M31 was -34332.168
$ws.Cells.Item(34, 8).Value = 234371.34  # H34 was 243730.3
$ws.Cells.Item(34, 9).Value = 32487.812  # I34 was 34627.168
$ws.Cells.Item(34, 11).Value = 32487.812  # K34 was 34627.168
$ws.Cells.Item(34, 13).Value = -32285.812  # M34 was -34425.168
$ws.Cells.Item(86, 8).Value = 2940.8  # H86 was 3067.5557
$ws.Cells.Item(86, 9).Value = 2200  # I86 was 2400
$ws.Cells.Item(86, 11).Value = 2200  # K86 was 2400
$ws.Cells.Item(86, 13).Value = -1077  # M86 was -1277
$ws.Cells.Item(89, 8).Value = 2940.8  # H89 was 3067.5557
$ws.Cells.Item(89, 9).Value = 2200  # I89 was 2400
$ws.Cells.Item(89, 11).Value = 11000  # K89 was 12000
$ws.Cells.Item(89, 13).Value = -5384  # M89 was -6384

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 5265392  # H4 was 6669192
$ws.Cells.Item(4, 9).Value = 30  # I4 was 1225
$ws.Cells.Item(4, 10).Value = 5884846.5  # J4 was 7695033
$ws.Cells.Item(4, 11).Value = 90  # K4 was 3675
$ws.Cells.Item(4, 12).Value = 17654539.5  # L4 was 23085099
$ws.Cells.Item(4, 13).Value = 22  # M4 was -3563
$ws.Cells.Item(4, 14).Value = -17654763.5  # N4 was -23085323
$ws.Cells.Item(14, 8).Value = 1367.375  # H14 was 1292.8235
$ws.Cells.Item(14, 9).Value = 1367.375  # I14 was 1292.8235
$ws.Cells.Item(14, 11).Value = 4102.125  # K14 was 3878.4705
$ws.Cells.Item(14, 13).Value = -3929.125  # M14 was -3705.4705
$ws.Cells.Item(108, 8).Value = 385.5  # H108 was 404.8
$ws.Cells.Item(108, 9).Value = 385.5  # I108 was 404.8
$ws.Cells.Item(108, 11).Value = 1156.5  # K108 was 1214.4
$ws.Cells.Item(108, 13).Value = 1723.5  # M108 was 1665.6
$ws.Cells.Item(113, 8).Value = 599.5161000000001  # H113 was 586.4828
$ws.Cells.Item(113, 9).Value = 563.35297  # I113 was 556.4286
$ws.Cells.Item(113, 10).Value = 643.4286  # J113 was 614.5333000000001
$ws.Cells.Item(113, 11).Value = 1690.05891  # K113 was 1669.2858
$ws.Cells.Item(113, 12).Value = 1930.2858  # L113 was 1843.5999
$ws.Cells.Item(113, 13).Value = 479.9410899999998  # M113 was 500.7142000000001
$ws.Cells.Item(113, 14).Value = -6270.2858  # N113 was -6183.5999
$ws.Cells.Item(131, 8).Value = 802.1539  # H131 was 786.0323
$ws.Cells.Item(131, 9).Value = 280.91666  # I131 was 295.54544
$ws.Cells.Item(131, 10).Value = 1033.8148  # J131 was 1055.8
$ws.Cells.Item(131, 11).Value = 842.7499799999999  # K131 was 886.63632
$ws.Cells.Item(131, 12).Value = 3101.4444  # L131 was 3167.4
$ws.Cells.Item(131, 13).Value = 4197.25002  # M131 was 4153.36368
$ws.Cells.Item(131, 14).Value = -13181.4444  # N131 was -13247.4

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8971.923000000001  # H80 was 3876.3635
$ws.Cells.Item(80, 9).Value = 13465  # I80 was 2000
$ws.Cells.Item(80, 10).Value = 7624  # J80 was 4064
$ws.Cells.Item(80, 11).Value = 13465  # K80 was 2000
$ws.Cells.Item(80, 12).Value = 7624  # L80 was 4064
$ws.Cells.Item(80, 13).Value = -12467  # M80 was -1002
$ws.Cells.Item(80, 14).Value = -9620  # N80 was -6060
$ws.Cells.Item(83, 8).Value = 8971.923000000001  # H83 was 3876.3635
$ws.Cells.Item(83, 9).Value = 13465  # I83 was 2000
$ws.Cells.Item(83, 10).Value = 7624  # J83 was 4064
$ws.Cells.Item(83, 11).Value = 67325  # K83 was 10000
$ws.Cells.Item(83, 12).Value = 38120  # L83 was 20320
$ws.Cells.Item(83, 13).Value = -62333  # M83 was -5008
$ws.Cells.Item(83, 14).Value = -48104  # N83 was -30304
$ws.Cells.Item(113, 8).Value = 1730.7273  # H113 was 2237.6667
$ws.Cells.Item(113, 9).Value = 1021.1429  # I113 was 1612.25
$ws.Cells.Item(113, 10).Value = 2972.5  # J113 was 2738
$ws.Cells.Item(113, 11).Value = 1021.1429  # K113 was 1612.25
$ws.Cells.Item(113, 12).Value = 2972.5  # L113 was 2738
$ws.Cells.Item(113, 13).Value = 1148.8571  # M113 was 557.75
$ws.Cells.Item(113, 14).Value = -7312.5  # N113 was -7078
$ws.Cells.Item(136, 8).Value = 16404.688  # H136 was 20812.416
$ws.Cells.Item(136, 10).Value = 16404.688  # J136 was 20812.416
$ws.Cells.Item(136, 12).Value = 49214.064  # L136 was 62437.24800000001
$ws.Cells.Item(136, 14).Value = -54314.064  # N136 was -67537.24800000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 646.08  # H22 was 642.38464
$ws.Cells.Item(22, 9).Value = 497.85715  # I22 was 474.66666
$ws.Cells.Item(22, 10).Value = 834.7273  # J22 was 871.0909
$ws.Cells.Item(22, 11).Value = 497.85715  # K22 was 474.66666
$ws.Cells.Item(22, 12).Value = 834.7273  # L22 was 871.0909
$ws.Cells.Item(22, 13).Value = -202.85715  # M22 was -179.66666
$ws.Cells.Item(22, 14).Value = -1424.7273  # N22 was -1461.0909
$ws.Cells.Item(27, 8).Value = 646.08  # H27 was 642.38464
$ws.Cells.Item(27, 9).Value = 497.85715  # I27 was 474.66666
$ws.Cells.Item(27, 10).Value = 834.7273  # J27 was 871.0909
$ws.Cells.Item(27, 11).Value = 497.85715  # K27 was 474.66666
$ws.Cells.Item(27, 12).Value = 834.7273  # L27 was 871.0909
$ws.Cells.Item(27, 13).Value = -390.85715  # M27 was -367.66666
$ws.Cells.Item(27, 14).Value = -1048.7273  # N27 was -1085.0909
$ws.Cells.Item(100, 8).Value = 1773.3125  # H100 was 1804.5625
$ws.Cells.Item(100, 9).Value = 1599.2222  # I100 was 1654.7778
$ws.Cells.Item(100, 11).Value = 1599.2222  # K100 was 1654.7778
$ws.Cells.Item(100, 13).Value = -1058.2222  # M100 was -1113.7778
$ws.Cells.Item(136, 8).Value = 47673.6  # H136 was 44253.49
$ws.Cells.Item(136, 9).Value = 32275.941  # I136 was 28781.309
$ws.Cells.Item(136, 10).Value = 95266.37  # J136 was 104595
$ws.Cells.Item(136, 11).Value = 96827.823  # K136 was 86343.927
$ws.Cells.Item(136, 12).Value = 285799.11  # L136 was 313785
$ws.Cells.Item(136, 13).Value = -94277.823  # M136 was -83793.927
$ws.Cells.Item(136, 14).Value = -290899.11  # N136 was -318885

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 112211.11  # H100 was 92027.55
$ws.Cells.Item(100, 9).Value = 84400.164  # I100 was 72485.86
$ws.Cells.Item(100, 10).Value = 167833  # J100 was 126225.5
$ws.Cells.Item(100, 11).Value = 168800.328  # K100 was 144971.72
$ws.Cells.Item(100, 12).Value = 335666  # L100 was 252451
$ws.Cells.Item(100, 13).Value = -168259.328  # M100 was -144430.72
$ws.Cells.Item(100, 14).Value = -336748  # N100 was -253533
$ws.Cells.Item(132, 8).Value = 36674.125  # H132 was 38027.055
$ws.Cells.Item(132, 9).Value = 22905.623  # I132 was 22943.623
$ws.Cells.Item(132, 10).Value = 92999.82000000001  # J132 was 113444.22
$ws.Cells.Item(132, 11).Value = 68716.86900000001  # K132 was 68830.86900000001
$ws.Cells.Item(132, 12).Value = 278999.46  # L132 was 340332.66
$ws.Cells.Item(132, 13).Value = -66186.86900000001  # M132 was -66300.86900000001
$ws.Cells.Item(132, 14).Value = -284059.46  # N132 was -345392.66
$ws.Cells.Item(136, 8).Value = 58254.223  # H136 was 51212.855
$ws.Cells.Item(136, 9).Value = 46915.547  # I136 was 38322.85
$ws.Cells.Item(136, 11).Value = 140746.641  # K136 was 114968.55
$ws.Cells.Item(136, 13).Value = -138196.641  # M136 was -112418.55

